# 2. normálforma, still ugly
$wb = $excel.ActiveWorkbook

# --- rename the existing sheet and add the new "2. normal form" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "1.normálforma"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2.normálforma"

# --- headers (same bold/centered style used for the header row on sheet 1) ---
$ws2.Range("A1").Value = "Rendszám"
$ws2.Range("B1").Value = "vkód"
$ws2.Range("C1").Value = "vevőnév"
$ws2.Range("D1").Value = "vevőcím"
$ws2.Range("E1").Value = "kelt"
$ws2.Range("F1").Value = "határidő"
$ws2.Range("G1").Value = "összérték"

$ws2.Range("J1").Value = "Rendszám"
$ws2.Range("K1").Value = "cikkszám"
$ws2.Range("L1").Value = "rendmenny"

$ws2.Range("O1").Value = "cikkszám"
$ws2.Range("P1").Value = "cikknév"
$ws2.Range("Q1").Value = "egysár"

$ws2.Range("A1:G1").Font.Bold = $true
$ws2.Range("A1:G1").HorizontalAlignment = -4108
$ws2.Range("J1:L1").Font.Bold = $true
$ws2.Range("J1:L1").HorizontalAlignment = -4108
$ws2.Range("O1:Q1").Font.Bold = $true
$ws2.Range("O1:Q1").HorizontalAlignment = -4108

# --- row 2 ---
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = "Vn1"
$ws2.Range("D2").Value = "VC1"
$ws2.Range("E2").Value = "09/14/2018"
$ws2.Range("F2").Value = "09/22/2018"
$ws2.Range("G2").Value = 550

$ws2.Range("J2").Value = 1
$ws2.Range("K2").Value = "CSZ1"
$ws2.Range("L2").Value = 2

$ws2.Range("O2").Value = "CSZ1"
$ws2.Range("P2").Value = "CN1"
$ws2.Range("Q2").Value = 125

# --- row 3 ---
$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = "Vn2"
$ws2.Range("D3").Value = "VC2"
$ws2.Range("E3").Value = "09/14/2018"
$ws2.Range("F3").Value = "09/22/2018"
$ws2.Range("G3").Value = 1375

$ws2.Range("J3").Value = 1
$ws2.Range("K3").Value = "CSZ2"
$ws2.Range("L3").Value = 1

$ws2.Range("O3").Value = "CSZ2"
$ws2.Range("P3").Value = "CN2"
$ws2.Range("Q3").Value = 300

$dateFmt = $ws1.Range("E2").NumberFormat
$ws2.Range("E2:F3").NumberFormat = $dateFmt

# --- row 4 (date cells present but empty, like the source sheet) ---
$ws2.Range("E4").NumberFormat = $dateFmt
$ws2.Range("F4").NumberFormat = $dateFmt

$ws2.Range("J4").Value = 2
$ws2.Range("K4").Value = "CSZ1"
$ws2.Range("L4").Value = 3

$ws2.Range("O4").Value = "CSZ3"
$ws2.Range("P4").Value = "CN3"
$ws2.Range("Q4").Value = 500

# --- row 5 ---
$ws2.Range("J5").Value = 2
$ws2.Range("K5").Value = "CSZ3"
$ws2.Range("L5").Value = 2

# --- column widths (mirrors sheet 1's "bestFit" columns, shifted into the new layout) ---
$ws2.Columns.Item(1).ColumnWidth = 9.17
$ws2.Columns.Item(5).ColumnWidth = 9.33
$ws2.Columns.Item(6).ColumnWidth = 9.33
$ws2.Columns.Item(8).ColumnWidth = 9.17
$ws2.Columns.Item(10).ColumnWidth = 10.5

# --- view state: sheet 1 keeps a (collapsed to single-area) selection, sheet 2 becomes active ---
$ws1.Select()
$ws1.Range("J1").Select()

$ws2.Select()
$ws2.Range("I4").Select()
$ws2.Application.ActiveWindow.Zoom = 100
$ws2.Application.ActiveWindow.ScrollColumn = 5
